# Fruta / hortaliza, semanal
# Insert two new weekly observations (date 44585) for Femacal de La Calera /
# Frutilla ahead of the existing row that used to be at row 225, shifting the
# rest of the "Frutilla" block down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 225-226; everything from the old row 225 onward
# shifts down to 227 onward, carrying its formatting (incl. the date style
# on column D) with it.
$ws.Range("A225:A226").EntireRow.Insert()

# New row 225: "Primera" quality observation dated 44585
$ws.Cells.Item(225, 1).Value = 3
$ws.Cells.Item(225, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(225, 3).Value = "Coquimbo"
$ws.Cells.Item(225, 4).Value = 44585
$ws.Cells.Item(225, 5).Value = 5
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100101
$ws.Cells.Item(225, 8).Value = "Berries"
$ws.Cells.Item(225, 9).Value = 100112025
$ws.Cells.Item(225, 10).Value = "Frutilla"
$ws.Cells.Item(225, 11).Value = "Sin especificar"
$ws.Cells.Item(225, 12).Value = "Primera"
$ws.Cells.Item(225, 13).Value = 120
$ws.Cells.Item(225, 14).Value = 5000
$ws.Cells.Item(225, 15).Value = 5000
$ws.Cells.Item(225, 16).Value = 5000
$ws.Cells.Item(225, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(225, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(225, 19).Value = 714
$ws.Cells.Item(225, 20).Value = 7

# New row 226: "Segunda" quality observation dated 44585
$ws.Cells.Item(226, 1).Value = 3
$ws.Cells.Item(226, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(226, 3).Value = "Coquimbo"
$ws.Cells.Item(226, 4).Value = 44585
$ws.Cells.Item(226, 5).Value = 5
$ws.Cells.Item(226, 6).Value = "Fruta"
$ws.Cells.Item(226, 7).Value = 100101
$ws.Cells.Item(226, 8).Value = "Berries"
$ws.Cells.Item(226, 9).Value = 100112025
$ws.Cells.Item(226, 10).Value = "Frutilla"
$ws.Cells.Item(226, 11).Value = "Sin especificar"
$ws.Cells.Item(226, 12).Value = "Segunda"
$ws.Cells.Item(226, 13).Value = 90
$ws.Cells.Item(226, 14).Value = 4000
$ws.Cells.Item(226, 15).Value = 4000
$ws.Cells.Item(226, 16).Value = 4000
$ws.Cells.Item(226, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(226, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(226, 19).Value = 571
$ws.Cells.Item(226, 20).Value = 7
